$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last three rows of data (rows 20, 21 and 22) first, while the
# rest of the table is still in its original position. (Two brand-new rows
# will be inserted below, for a net change of -1 row.)
$ws.Rows.Item(20).Resize(3).Delete()

# Insert two new rows at the top of the data (row 2), pushing the
# remaining rows down by two.
$ws.Rows.Item(2).Resize(2).Insert()

# The inserted rows pick up the header row's formatting by default;
# strip it so the new data rows look like the other plain data rows.
$ws.Range("A2:C3").ClearFormats()

# Populate the two newly inserted rows with the new accelerometer readings.
$ws.Range("A2").Value = -1.921436786651612
$ws.Range("B2").Value = 1.638089656829834
$ws.Range("C2").Value = 0.3902863562107086

$ws.Range("A3").Value = -1.84520435333252
$ws.Range("B3").Value = 1.668948650360107
$ws.Range("C3").Value = 0.1994338035583496
